$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "nik" -> "nip/nrpk"
$ws.Range("A1").Value = "nip/nrpk"

# kode_akses column (G) - normalize all rows to 12345
$ws.Range("G3").Value = 12345
$ws.Range("G4").Value = 12345
$ws.Range("G5").Value = 12345
$ws.Range("G6").Value = 12345

# pph (E) and diterima (F) become formulas
$ws.Range("E2").Formula = "=2%*D2"
$ws.Range("F2").Formula = "=D2-E2"

$ws.Range("E3:E6").Formula = "=2%*D3"
$ws.Range("F3:F6").Formula = "=D3-E3"

# Remove the warning column H entirely
$ws.Columns("H").Delete()

# Update the active selection shown when the sheet was saved
$null = $ws.Range("E2").Select()
